# Add "other relevant files" sample data (Category/Type/Item/Size) below the
# existing header row on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Clothing", "T-Shirt", "Polo",      "Small"),
    @("Clothing", "T-Shirt", "Polo",      "Medium"),
    @("Clothing", "T-Shirt", "Crew Neck", "Large"),
    @("Clothing", "T-Shirt", "Crew Neck", "Extra large")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Match the author's final cursor position/selection.
$ws.Range("B12").Select()
